$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in header (A1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 22:51"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4690274
$ws.Range("C4").Value = 55289
$ws.Range("D4").Value = 2305798
$ws.Range("E4").Value = 2228211
$ws.Range("G4").Value = 980
$ws.Range("H4").Value = 156265

# Row 8 - Sudafrica
$ws.Range("B8").Value = 493183
$ws.Range("C8").Value = 11014
$ws.Range("D8").Value = 326171
$ws.Range("E8").Value = 159007
$ws.Range("G8").Value = 193
$ws.Range("H8").Value = 8005

# Row 21 - Alemania
$ws.Range("B21").Value = 210562
$ws.Range("C21").Value = 909
$ws.Range("E21").Value = 9038
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 9224

# Row 31 - Ecuador
$ws.Range("E31").Value = 43609
$ws.Range("G31").Value = 45
$ws.Range("H31").Value = 5702

# Row 48 - Guatemala
$ws.Range("B48").Value = 49789
$ws.Range("C48").Value = 963
$ws.Range("D48").Value = 36816
$ws.Range("E48").Value = 11049
$ws.Range("G48").Value = 57
$ws.Range("H48").Value = 1924

# Row 56/57 - Ghana overtakes Suiza (swap order + update Ghana's data)
$ws.Range("A56").Value = "Ghana"
$ws.Range("B56").Value = 35501
$ws.Range("C56").Value = 359
$ws.Range("D56").Value = 32096
$ws.Range("E56").Value = 3223
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 182

$ws.Range("A57").Value = "Suiza"
$ws.Range("B57").Value = 35232
$ws.Range("C57").Value = 210
$ws.Range("D57").Value = 31100
$ws.Range("E57").Value = 2151
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 1981

# Row 70/71 - Costa Rica overtakes Etiopia (swap order + update Costa Rica's data)
$ws.Range("A70").Value = "Costa Rica"
$ws.Range("B70").Value = 17820
$ws.Range("C70").Value = 530
$ws.Range("D70").Value = 4404
$ws.Range("E70").Value = 13266
$ws.Range("G70").Value = 10
$ws.Range("H70").Value = 150

$ws.Range("A71").Value = "Etiopia"
$ws.Range("B71").Value = 17530
$ws.Range("C71").Value = 915
$ws.Range("D71").Value = 6950
$ws.Range("E71").Value = 10306
$ws.Range("G71").Value = 11
$ws.Range("H71").Value = 274

# Row 94 - Guinea
$ws.Range("B94").Value = 7308
$ws.Range("C94").Value = 66
$ws.Range("D94").Value = 6458
$ws.Range("E94").Value = 804

# Row 109 - Libia
$ws.Range("B109").Value = 3621
$ws.Range("C109").Value = 183
$ws.Range("D109").Value = 618
$ws.Range("E109").Value = 2929
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 74
